$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item(5)

# Insert a new "category" column before the existing "date" column (col I -> shifts to J)
$ws5.Columns.Item(9).Insert()
$ws5.Cells.Item(1,9).Value = "category"
$ws5.Cells.Item(2,9).Value = "normal"
$ws5.Cells.Item(3,9).Value = "normal"

# Append two new trailing columns: source_file, index (after legislator_id, now col L)
$ws5.Cells.Item(1,13).Value = "source_file"
$ws5.Cells.Item(2,13).Value = "tmp55951"
$ws5.Cells.Item(3,13).Value = "tmp55951"

$ws5.Cells.Item(1,14).Value = "index"
$ws5.Cells.Item(2,14).Value = 61
$ws5.Cells.Item(3,14).Value = 62

# Copy the bold/bordered header style onto the two new trailing header cells
$ws5.Range("L1").Copy() | Out-Null
$ws5.Range("M1:N1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
